$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 404.53333
$ws.Range("I33").Value = 391.70834
$ws.Range("K33").Value = 391.70834
$ws.Range("M33").Value = -162.70834

$ws.Range("H40").Value = 2072.4614
$ws.Range("I40").Value = 1980.25
$ws.Range("J40").Value = 2220
$ws.Range("K40").Value = 1980.25
$ws.Range("L40").Value = 2220
$ws.Range("M40").Value = -1805.25
$ws.Range("N40").Value = -2570

$ws.Range("H74").Value = 3314.7144
$ws.Range("I74").Value = 2950.75
$ws.Range("K74").Value = 2950.75
$ws.Range("M74").Value = -2014.75

$ws.Range("H77").Value = 3314.7144
$ws.Range("I77").Value = 2950.75
$ws.Range("K77").Value = 14753.75
$ws.Range("M77").Value = -10073.75

$ws.Range("H88").Value = 1124377.4
$ws.Range("I88").Value = 1786.6
$ws.Range("J88").Value = 2059869.6
$ws.Range("K88").Value = 1786.6
$ws.Range("L88").Value = 2059869.6
$ws.Range("M88").Value = -1380.6
$ws.Range("N88").Value = -2060681.6

$ws.Range("H91").Value = 1124377.4
$ws.Range("I91").Value = 1786.6
$ws.Range("J91").Value = 2059869.6
$ws.Range("K91").Value = 1786.6
$ws.Range("L91").Value = 2059869.6
$ws.Range("M91").Value = -382.5999999999999
$ws.Range("N91").Value = -2062677.6

$ws.Range("H138").Value = 792144.8
$ws.Range("I138").Value = 1393.8125
$ws.Range("J138").Value = 1116555.5
$ws.Range("K138").Value = 4181.4375
$ws.Range("L138").Value = 3349666.5
$ws.Range("M138").Value = 958.5625
$ws.Range("N138").Value = -3359946.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 674.9048
$ws.Range("I2").Value = 475.2143
$ws.Range("J2").Value = 1074.2858
$ws.Range("K2").Value = 475.2143
$ws.Range("L2").Value = 1074.2858
$ws.Range("M2").Value = -362.2143
$ws.Range("N2").Value = -1300.2858

$ws.Range("H32").Value = 8386.267
$ws.Range("I32").Value = 6950.8276
$ws.Range("J32").Value = 50014
$ws.Range("K32").Value = 6950.8276
$ws.Range("L32").Value = 50014
$ws.Range("M32").Value = -6663.8276
$ws.Range("N32").Value = -50588

$ws.Range("H97").Value = 690.6429000000001
$ws.Range("I97").Value = 423.45456
$ws.Range("K97").Value = 423.45456
$ws.Range("M97").Value = 72.54543999999999

$ws.Range("H114").Value = 29499.5
$ws.Range("J114").Value = 29499.5
$ws.Range("L114").Value = 29499.5
$ws.Range("N114").Value = -38177.5

$ws.Range("H116").Value = 674.9048
$ws.Range("I116").Value = 475.2143
$ws.Range("J116").Value = 1074.2858
$ws.Range("K116").Value = 475.2143
$ws.Range("L116").Value = 1074.2858
$ws.Range("M116").Value = 1818.7857
$ws.Range("N116").Value = -5662.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 27003.334
$ws.Range("I57").Value = 26000
$ws.Range("J57").Value = 27505
$ws.Range("K57").Value = 26000
$ws.Range("L57").Value = 27505
$ws.Range("M57").Value = -25280
$ws.Range("N57").Value = -28945

$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 1282.6154
$ws.Range("I107").Value = 923.7273
$ws.Range("K107").Value = 923.7273
$ws.Range("M107").Value = 996.2727

$ws.Range("H132").Value = 10000000
$ws.Range("J132").Value = 10000000
$ws.Range("L132").Value = 10000000
$ws.Range("N132").Value = -10010120

$ws.Range("H136").Value = 27003.334
$ws.Range("I136").Value = 26000
$ws.Range("J136").Value = 27505
$ws.Range("K136").Value = 26000
$ws.Range("L136").Value = 27505
$ws.Range("M136").Value = -20900
$ws.Range("N136").Value = -37705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 9093150
$ws.Range("J62").Value = 200000000
$ws.Range("L62").Value = 200000000
$ws.Range("N62").Value = -200001248

$ws.Range("H65").Value = 9093150
$ws.Range("J65").Value = 200000000
$ws.Range("L65").Value = 1000000000
$ws.Range("N65").Value = -1000006240

$ws.Range("I134").Value = 1689.75
$ws.Range("K134").Value = 5069.25
$ws.Range("M134").Value = -2534.25

$ws.Range("H141").Value = 1512200
$ws.Range("J141").Value = 1512200
$ws.Range("L141").Value = 1512200
$ws.Range("N141").Value = -1522560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 52631710
$ws.Range("I12").Value = 200000210
$ws.Range("J12").Value = 106
$ws.Range("K12").Value = 600000630
$ws.Range("L12").Value = 318
$ws.Range("M12").Value = -600000457
$ws.Range("N12").Value = -664

$ws.Range("H113").Value = 667.85187
$ws.Range("I113").Value = 577.5
$ws.Range("J113").Value = 683.56525
$ws.Range("K113").Value = 1732.5
$ws.Range("L113").Value = 2050.69575
$ws.Range("M113").Value = 437.5
$ws.Range("N113").Value = -6390.69575

$ws.Range("H117").Value = 1326.1
$ws.Range("I117").Value = 625.6667
$ws.Range("J117").Value = 1626.2858
$ws.Range("K117").Value = 1877.0001
$ws.Range("L117").Value = 4878.857400000001
$ws.Range("M117").Value = 1564.9999
$ws.Range("N117").Value = -11762.8574

$ws.Range("H129").Value = 16027060
$ws.Range("J129").Value = 4387640
$ws.Range("L129").Value = 13162920
$ws.Range("N129").Value = -13172920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2711.074
$ws.Range("I40").Value = 1827.0454
$ws.Range("J40").Value = 6600.8
$ws.Range("K40").Value = 1827.0454
$ws.Range("L40").Value = 6600.8
$ws.Range("M40").Value = -1691.0454
$ws.Range("N40").Value = -6872.8

$ws.Range("H46").Value = 3764.6428
$ws.Range("I46").Value = 554.5
$ws.Range("K46").Value = 554.5
$ws.Range("M46").Value = -366.5

$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 35000
$ws.Range("N134").Value = -45140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
